$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeGetResults")

$ws.Range("A2").Value = "https://stackoverflow.com/questions/46854451/pip-install-r-requirements-txt-errno-2-no-such-file-or-directory-requiremen"
$ws.Range("B2").Value = 200
